# Auto-generated script to apply numeric corrections to Carbuncle_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7579.4
$ws.Range("I19").Value = 33940.332
$ws.Range("J19").Value = 989.1667
$ws.Range("K19").Value = 33940.332
$ws.Range("L19").Value = 989.1667
$ws.Range("M19").Value = -33765.332
$ws.Range("N19").Value = -1339.1667
$ws.Range("H74").Value = 4544.778
$ws.Range("J74").Value = 4487.5
$ws.Range("L74").Value = 4487.5
$ws.Range("N74").Value = -6359.5
$ws.Range("H76").Value = 3344.492
$ws.Range("I76").Value = 3117.2456
$ws.Range("K76").Value = 3117.2456
$ws.Range("M76").Value = -2802.2456
$ws.Range("H77").Value = 4544.778
$ws.Range("J77").Value = 4487.5
$ws.Range("L77").Value = 22437.5
$ws.Range("N77").Value = -31797.5
$ws.Range("H79").Value = 3344.492
$ws.Range("I79").Value = 3117.2456
$ws.Range("K79").Value = 3117.2456
$ws.Range("M79").Value = -2025.2456
$ws.Range("H107").Value = 822.6
$ws.Range("I107").Value = 1213.3334
$ws.Range("J107").Value = 236.5
$ws.Range("K107").Value = 1213.3334
$ws.Range("L107").Value = 236.5
$ws.Range("M107").Value = 706.6666
$ws.Range("N107").Value = -4076.5
$ws.Range("H132").Value = 911.6923
$ws.Range("I132").Value = 527.7222
$ws.Range("J132").Value = 1775.625
$ws.Range("K132").Value = 1583.1666
$ws.Range("L132").Value = 5326.875
$ws.Range("M132").Value = 946.8334
$ws.Range("N132").Value = -10386.875
$ws.Range("H137").Value = 1786.7222
$ws.Range("I137").Value = 1837.4
$ws.Range("J137").Value = 1533.3334
$ws.Range("K137").Value = 5512.200000000001
$ws.Range("L137").Value = 4600.0002
$ws.Range("M137").Value = -2962.200000000001
$ws.Range("N137").Value = -9700.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 13394.25
$ws.Range("J43").Value = 13059
$ws.Range("L43").Value = 13059
$ws.Range("N43").Value = -13685
$ws.Range("H61").Value = 4034.7742
$ws.Range("I61").Value = 2649.7273
$ws.Range("J61").Value = 7420.4443
$ws.Range("K61").Value = 2649.7273
$ws.Range("L61").Value = 7420.4443
$ws.Range("M61").Value = -2437.7273
$ws.Range("N61").Value = -7844.4443
$ws.Range("H132").Value = 2933.9062
$ws.Range("I132").Value = 1644.2273
$ws.Range("J132").Value = 5771.2
$ws.Range("K132").Value = 4932.6819
$ws.Range("L132").Value = 17313.6
$ws.Range("M132").Value = -2402.6819
$ws.Range("N132").Value = -22373.6
$ws.Range("H136").Value = 4034.7742
$ws.Range("I136").Value = 2649.7273
$ws.Range("J136").Value = 7420.4443
$ws.Range("K136").Value = 7949.1819
$ws.Range("L136").Value = 22261.3329
$ws.Range("M136").Value = -5399.1819
$ws.Range("N136").Value = -27361.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1336.6666
$ws.Range("I105").Value = 1005
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1005
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 742
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3442.956
$ws.Range("I31").Value = 776.375
$ws.Range("J31").Value = 4889.2373
$ws.Range("K31").Value = 776.375
$ws.Range("L31").Value = 4889.2373
$ws.Range("M31").Value = -481.375
$ws.Range("N31").Value = -5479.2373
$ws.Range("H34").Value = 3442.956
$ws.Range("I34").Value = 776.375
$ws.Range("J34").Value = 4889.2373
$ws.Range("K34").Value = 776.375
$ws.Range("L34").Value = 4889.2373
$ws.Range("M34").Value = -574.375
$ws.Range("N34").Value = -5293.2373
$ws.Range("H132").Value = 3130.25
$ws.Range("I132").Value = 2955.4167
$ws.Range("K132").Value = 8866.250100000001
$ws.Range("M132").Value = -6336.250100000001
$ws.Range("H134").Value = 3416.4546
$ws.Range("I134").Value = 4434.8076
$ws.Range("J134").Value = 1945.5
$ws.Range("K134").Value = 13304.4228
$ws.Range("L134").Value = 5836.5
$ws.Range("M134").Value = -10769.4228
$ws.Range("N134").Value = -10906.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 345097.28
$ws.Range("I5").Value = 794.125
$ws.Range("J5").Value = 409153.7
$ws.Range("K5").Value = 2382.375
$ws.Range("L5").Value = 1227461.1
$ws.Range("M5").Value = -2270.375
$ws.Range("N5").Value = -1227685.1
$ws.Range("H68").Value = 2191668.2
$ws.Range("I68").Value = 9736446
$ws.Range("J68").Value = 1248.7742
$ws.Range("K68").Value = 29209338
$ws.Range("L68").Value = 3746.3226
$ws.Range("M68").Value = -29208527
$ws.Range("N68").Value = -5368.3226
$ws.Range("H71").Value = 2191668.2
$ws.Range("I71").Value = 9736446
$ws.Range("J71").Value = 1248.7742
$ws.Range("K71").Value = 87628014
$ws.Range("L71").Value = 11238.9678
$ws.Range("M71").Value = -87623958
$ws.Range("N71").Value = -19350.9678
$ws.Range("H131").Value = 1515.1978
$ws.Range("I131").Value = 772.5
$ws.Range("J131").Value = 1549.3448
$ws.Range("K131").Value = 2317.5
$ws.Range("L131").Value = 4648.0344
$ws.Range("M131").Value = 2722.5
$ws.Range("N131").Value = -14728.0344
$ws.Range("H135").Value = 345097.28
$ws.Range("I135").Value = 794.125
$ws.Range("J135").Value = 409153.7
$ws.Range("K135").Value = 7147.125
$ws.Range("L135").Value = 3682383.3
$ws.Range("M135").Value = -4612.125
$ws.Range("N135").Value = -3687453.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 48057.637
$ws.Range("I7").Value = 58176
$ws.Range("K7").Value = 58176
$ws.Range("M7").Value = -58064
$ws.Range("H126").Value = 48057.637
$ws.Range("I126").Value = 58176
$ws.Range("K126").Value = 174528
$ws.Range("M126").Value = -172058
$ws.Range("H136").Value = 8335506.5
$ws.Range("I136").Value = 829.6875
$ws.Range("J136").Value = 13891958
$ws.Range("K136").Value = 2489.0625
$ws.Range("L136").Value = 41675874
$ws.Range("M136").Value = 60.9375
$ws.Range("N136").Value = -41680974

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 811.5714
$ws.Range("I107").Value = 811.5714
$ws.Range("K107").Value = 2434.7142
$ws.Range("M107").Value = -514.7142000000003
$ws.Range("H132").Value = 2344.889
$ws.Range("I132").Value = 1852.8422
$ws.Range("J132").Value = 2894.8235
$ws.Range("K132").Value = 5558.5266
$ws.Range("L132").Value = 8684.470499999999
$ws.Range("M132").Value = -3028.5266
$ws.Range("N132").Value = -13744.4705
$ws.Range("H136").Value = 4585
$ws.Range("I136").Value = 833.1
$ws.Range("K136").Value = 2499.3
$ws.Range("M136").Value = 50.69999999999982
